# Update the roster table (A2:C19) on the active sheet so that each
# player row reflects the corrected Position / Team assignment, adds
# "Naz Reid" as a new player, and removes "Julian Champagnie".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Russell Westbrook", "PG",    "Denver Nuggets"),
    @("Jalen Suggs",        "PG,SG", "Orlando Magic"),
    @("Dejounte Murray",    "PG,SG", "New Orleans Pelicans"),
    @("Chris Paul",         "PG",    "San Antonio Spurs"),
    @("Pascal Siakam",      "SF,PF", "Indiana Pacers"),
    @("Deni Avdija",        "SF,PF", "Portland Trail Blazers"),
    @("Jerami Grant",       "SF,PF", "Portland Trail Blazers"),
    @("Naz Reid",           "PF,C",  "Minnesota Timberwolves"),
    @("Nikola Jokic",       "C",     "Denver Nuggets"),
    @("Rudy Gobert",        "C",     "Minnesota Timberwolves"),
    @("Clint Capela",       "C",     "Atlanta Hawks"),
    @("Jalen Green",        "PG,SG", "Houston Rockets"),
    @("Jaylen Brown",       "SG,SF", "Boston Celtics"),
    @("Jakob Poeltl",       "C",     "Toronto Raptors"),
    @("Ayo Dosunmu",        "SG,SF", "Chicago Bulls"),
    @("Paolo Banchero",     "SF,PF", "Orlando Magic"),
    @("Chet Holmgren",      "PF,C",  "Oklahoma City Thunder"),
    @("Bogdan Bogdanovic",  "SG,SF", "Atlanta Hawks")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
